# Hortaliza, Agrícola del Norte S.A. de Arica - Cebolla
# Weekly fruit/vegetable price update: insert 8 new daily observations
# (rows 1096-1103) ahead of the existing history block, pushing the
# previously-existing rows 1096-1133 down to 1104-1141.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 blank rows at row 1096 (each Insert() pushes everything at/after
# the target row down by one, so doing this 8 times at the same index
# opens up exactly 8 fresh rows at 1096..1103).
for ($i = 0; $i -lt 8; $i++) {
    $ws.Rows.Item(1096).Insert()
}

# Columns that are constant across this whole sub-block of data.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$catId     = 100112004
$categoria = "Cebolla"
$unidad    = "`$/malla 18 kilos"
$kgUnid    = 18
$clasif    = "Hortaliza"

# New data for rows 1096-1103 (Fecha, Variedad, Calidad, Volumen, PrecioMin,
# PrecioMax, PrecioPromedio, Origen, Precio $/Kg).
$newRows = @(
    @{ Row=1096; Fecha=44939; Variedad="Morada(o)";        Calidad="1a (cosecha)"; Vol=400; Min=6000; Max=7000; Prom=6500; Origen="Región de Arica y Parinacota"; PrecioKg=361 },
    @{ Row=1097; Fecha=44939; Variedad="Morada(o)";        Calidad="2a (cosecha)"; Vol=500; Min=5000; Max=6000; Prom=5500; Origen="Región de Arica y Parinacota"; PrecioKg=306 },
    @{ Row=1098; Fecha=44939; Variedad="Morada(o)";        Calidad="3a (cosecha)"; Vol=600; Min=4000; Max=5000; Prom=4500; Origen="Región de Arica y Parinacota"; PrecioKg=250 },
    @{ Row=1099; Fecha=44939; Variedad="Sin especificar";  Calidad="1a (cosecha)"; Vol=600; Min=8000; Max=9000; Prom=8583; Origen="Región de Arica y Parinacota"; PrecioKg=477 },
    @{ Row=1100; Fecha=44939; Variedad="Sin especificar";  Calidad="2a (cosecha)"; Vol=500; Min=6000; Max=7000; Prom=6600; Origen="Región de Arica y Parinacota"; PrecioKg=367 },
    @{ Row=1101; Fecha=44939; Variedad="Sin especificar";  Calidad="3a (cosecha)"; Vol=300; Min=4000; Max=5000; Prom=4667; Origen="Región de Arica y Parinacota"; PrecioKg=259 },
    @{ Row=1102; Fecha=44939; Variedad="Sin especificar";  Calidad="Primera";      Vol=450; Min=6000; Max=7000; Prom=6500; Origen="Perú";                        PrecioKg=361 },
    @{ Row=1103; Fecha=44939; Variedad="Sin especificar";  Calidad="Segunda";      Vol=500; Min=5000; Max=6000; Prom=5500; Origen="Perú";                        PrecioKg=306 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $r.Variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Vol
    $ws.Cells.Item($row, 11).Value = $r.Min
    $ws.Cells.Item($row, 12).Value = $r.Max
    $ws.Cells.Item($row, 13).Value = $r.Prom
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PrecioKg
    $ws.Cells.Item($row, 17).Value = $kgUnid
    $ws.Cells.Item($row, 18).Value = $clasif
}
